$d = $word.ActiveDocument

# --- Paragraph 1: "<student group>" placeholder -> MERGEFIELD student.group / student.name ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("`${group} `${student}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "paragraph 1 anchor text not found" }
$p1 = $r1.Paragraphs(1)
$pr1 = $p1.Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="0D9848D2" w14:textId="77777777" w:rsidR="00B7089A" w:rsidRPr="00E87528" w:rsidRDefault="00B56F5B" w:rsidP="001F356E" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:i/><w:iCs/><w:color w:val="000000"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r w:rsidRPr="00297955"><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> Студенту </w:t></w:r><w:r w:rsidRPr="00297955"><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="00E87528"><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>MERGEFIELD</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">  ${</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>student</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText>.</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>group</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">}  \* </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>MERGEFORMAT</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/></w:rPr><w:t>student</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/></w:rPr><w:t>group</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>MERGEFIELD</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">  ${</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>student</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText>.</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>name</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">}  \* </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:instrText>MERGEFORMAT</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/></w:rPr><w:t>student</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:lang w:val="ru-RU"/></w:rPr><w:t>}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>
'@
$pr1.InsertXML($xml1)

# --- Paragraph 2: project theme placeholder -> MERGEFIELD topic.title ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("`${projectTheme} `"", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "paragraph 2 anchor text not found" }
$p2 = $r2.Paragraphs(1)
$pr2 = $p2.Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="59E31CCA" w14:textId="5EBF36C5" w:rsidR="00B56F5B" w:rsidRDefault="00B56F5B" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">1. Тема </w:t></w:r><w:bookmarkStart w:id="0" w:name="_Hlk197954814"/><w:r w:rsidR="00AD0F3D"><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>работы</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">  "</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:instrText>MERGEFIELD</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">  ${</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:instrText>topic</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText>.</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:instrText>title</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve">}  \* </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:instrText>MERGEFORMAT</w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t>«${</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:u w:val="single"/></w:rPr><w:t>topic</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:u w:val="single"/></w:rPr><w:t>title</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t>}»</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:u w:val="single"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">                                </w:t></w:r></w:p>
'@
$pr2.InsertXML($xml2)

Write-Host "edit complete"
